# Update cryptocurrency price/volume data on the active sheet
# Prices are stored as text (some use apostrophe-prefix to force text
# for numeric-looking strings, matching how Excel treats quote-prefixed entries),
# and volume/percent-change values are plain text with surrounding spaces.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.884.49"
$ws.Range("E2").Value = "  +2.77%  "
$ws.Range("D3").Value = "1.859.40"
$ws.Range("E3").Value = "  +2.10%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'246.45"
$ws.Range("E5").Value = "  +2.13%  "
$ws.Range("D6").Value = "'0.6373"
$ws.Range("E6").Value = "  +3.65%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "'0.3000"
$ws.Range("E8").Value = "  +4.00%  "
$ws.Range("D9").Value = "'0.07472"
$ws.Range("E9").Value = "  +1.96%  "
$ws.Range("D10").Value = "'24.56"
$ws.Range("E10").Value = "  +7.25%  "
$ws.Range("D11").Value = "'0.07680"
$ws.Range("E11").Value = "  +0.35%  "
$ws.Range("D12").Value = "1.866.30"
$ws.Range("E12").Value = "  +2.46%  "
$ws.Range("D13").Value = "'5.052"
$ws.Range("E13").Value = "  +2.27%  "
$ws.Range("D14").Value = "'0.6910"
$ws.Range("E14").Value = "  +4.89%  "
$ws.Range("D15").Value = "'84.33"
$ws.Range("E15").Value = "  +3.30%  "
$ws.Range("D16").Value = "'0.000009357"
$ws.Range("E16").Value = "  +4.04%  "
$ws.Range("D17").Value = "'6.089"
$ws.Range("E17").Value = "  +4.56%  "
$ws.Range("D18").Value = "29.851.04"
$ws.Range("E18").Value = "  +2.79%  "
$ws.Range("D19").Value = "2.114.13"
$ws.Range("E19").Value = "  +2.72%  "
$ws.Range("D20").Value = "'238.40"
$ws.Range("E20").Value = "  +0.62%  "
$ws.Range("E21").Value = "  +1.94%  "
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").Value = "'7.352"
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").Value = "'159.26"
$ws.Range("E25").Value = "  +1.28%  "
$ws.Range("D26").Value = "'0.1418"
$ws.Range("E26").Value = "  +1.29%  "
$ws.Range("D27").Value = "'8.580"
$ws.Range("E27").Value = "  +1.93%  "
$ws.Range("D28").Value = "'17.96"
$ws.Range("E28").Value = "  +2.30%  "
$ws.Range("D29").Value = "'1.506"
$ws.Range("E29").Value = "  +1.66%  "
$ws.Range("D30").Value = "'0.06067"
$ws.Range("E30").Value = "  +9.21%  "
$ws.Range("D31").Value = "'1.283"
$ws.Range("E31").Value = "  +6.49%  "
$ws.Range("D32").Value = "'4.135"
$ws.Range("E32").Value = "  +1.06%  "
$ws.Range("D33").Value = "'4.141"
$ws.Range("E33").Value = "  +1.34%  "
$ws.Range("D34").Value = "'1.894"
$ws.Range("E34").Value = "  +4.60%  "
$ws.Range("D35").Value = "'1.166"
$ws.Range("E35").Value = "  +3.22%  "
$ws.Range("E36").Value = "  -0.62%  "
$ws.Range("D37").Value = "'2.607"
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").Value = "'2.863"
$ws.Range("E38").Value = "  +1.40%  "
$ws.Range("D39").Value = "'0.01797"
$ws.Range("E39").Value = "  +2.62%  "
$ws.Range("D40").Value = "1.223.26"
$ws.Range("E40").Value = "  +1.33%  "
$ws.Range("D41").Value = "'0.9312"
$ws.Range("E41").Value = "  +4.56%  "
$ws.Range("D42").Value = "'6.301"
$ws.Range("E42").Value = "  -0.61%  "
$ws.Range("D43").Value = "'1.000"
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("D44").Value = "2.034.45"
$ws.Range("E44").Value = "  +3.72%  "
$ws.Range("D45").Value = "'102.42"
$ws.Range("E45").Value = "  +1.44%  "
$ws.Range("D46").Value = "'66.29"
$ws.Range("E46").Value = "  +2.95%  "
$ws.Range("E47").Value = "  +4.96%  "
$ws.Range("D48").Value = "'0.5094"
$ws.Range("E48").Value = "  +0.29%  "
$ws.Range("D49").Value = "'9.310"
$ws.Range("E49").Value = "  +3.41%  "
$ws.Range("D50").Value = "'0.4090"
$ws.Range("E50").Value = "  +2.46%  "
$ws.Range("D51").Value = "'0.1144"
$ws.Range("E51").Value = "  +3.27%  "
